$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Insert a brand-new paragraph right after the "strike-through / small caps"
# paragraph (paragraph 2) and before the "subscripts / superscript" paragraph,
# mirroring the existing paragraphs' own plain (un-styled) run boundaries.
$strikeThroughPara = $tr.Paragraphs(2, 1)
$null = $strikeThroughPara.InsertAfter([char]13 + "Here is some underlined text.")

# The inserted text now forms its own paragraph (paragraph 3).
$newPara = $tr.Paragraphs(3, 1)

# "Here is " -> no extra formatting
# "some "    -> underline
# "underlined" -> italic + underline
# " "        -> underline
# "text"     -> bold + underline
# "."        -> no extra formatting
$newPara.Characters(9, 5).Font.Underline = $true
$newPara.Characters(14, 10).Font.Underline = $true
$newPara.Characters(14, 10).Font.Italic = $true
$newPara.Characters(24, 1).Font.Underline = $true
$newPara.Characters(25, 4).Font.Underline = $true
$newPara.Characters(25, 4).Font.Bold = $true
